{"js": "// Change the year in the astromap link: .../GaNight/2018/ -> .../GaNight/2022/\n// The sentence (\"por Jenik Hollan, CzechGlobe (http://.../GaNight/2018/).\") is\n// split across many differently-formatted runs (with proofErr spell-check\n// wrappers in between) in the original document. We locate the paragraph\n// that holds this sentence, compute its updated plain text, and rewrite the\n// paragraph's content with that single corrected string.\n\nconst body = context.document.body;\n\n// \"CzechGlobe\" only occurs once in the document, inside the sentence we need\n// to edit, so it is a safe, unique anchor for locating the paragraph.\nconst results = body.search(\"CzechGlobe\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the 'CzechGlobe' astromap credit line.\");\n}\n\nconst hit = results.items[0];\nconst paragraphs = hit.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nparagraph.load(\"text\");\nawait context.sync();\n\nconst oldText = paragraph.text;\nconst newText = oldText.replace(\"2018\", \"2022\");\n\nif (newText !== oldText) {\n  // Insert the corrected text as a brand-new (plain) run right before the\n  // paragraph's existing content, then delete everything that follows it,\n  // leaving only the corrected sentence in the paragraph.\n  const startRange = paragraph.getRange(\"Start\");\n  const insertedRange = startRange.insertText(newText, Word.InsertLocation.before);\n  await context.sync();\n\n  const afterInserted = insertedRange.getRange(Word.RangeLocation.after);\n  const paragraphEnd = paragraph.getRange(\"End\");\n  const oldContentRange = afterInserted.expandTo(paragraphEnd);\n  oldContentRange.delete();\n  await context.sync();\n}\n", "ps1": "# Change the year in the astromap link: .../GaNight/2018/ -> .../GaNight/2022/\n# The sentence (\"por Jenik Hollan, CzechGlobe (http://.../GaNight/2018/).\") is\n# split across many differently-formatted runs (with spell-check proofErr\n# wrappers in between) in the original document. We locate the paragraph\n# that holds this sentence, compute its updated plain text, and rewrite the\n# paragraph's content with that single corrected string.\n\n$d = $word.ActiveDocument\n\n# \"CzechGlobe\" only occurs once in the document, inside the sentence we need\n# to edit, so it is a safe, unique anchor for locating the paragraph.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"CzechGlobe\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the 'CzechGlobe' astromap credit line.\"\n}\n\n# Grow the found hit out to the whole enclosing paragraph (wdParagraph),\n# which also pulls in the trailing paragraph mark.\n$rng.Expand(4) | Out-Null\n\n$origStart = $rng.Start\n$origEnd = $rng.End\n\n$oldText = $rng.Text.TrimEnd(\"`r\")\n$newText = $oldText -replace \"2018\", \"2022\"\n\nif ($newText -ne $oldText) {\n    # Insert the corrected sentence as a brand-new (plain) run right before\n    # the paragraph's existing content.\n    $startR = $rng.Duplicate\n    $startR.Collapse(1)  # wdCollapseStart\n    $startR.InsertBefore($newText)\n\n    # Delete everything that followed the original insertion point (i.e. all\n    # of the old, differently-formatted runs), shifted right by the length\n    # of the text we just inserted, leaving only the corrected sentence.\n    $shift = $newText.Length\n    $tailR = $d.Range($origStart + $shift, $origEnd + $shift)\n    $tailR.Delete()\n}\n"}
